$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.022.65"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.902.08"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7405"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3064"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06899"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08042"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7621"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("D13").Value = "1.915.54"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.229"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "30.030.33"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.070"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007752"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "2.152.00"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.064"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.294"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1261"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.035"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.353"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.534"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.276"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.037"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05434"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.293"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7349"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01942"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.793"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.268"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4439"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.954"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8325"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.620"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.762"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").Value = "2.054.48"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.46%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1170"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.29%  "
